# Commit: "Any user can vote anonymously in a poll and the vote is recorded
# correctly, moving to Chores"
#
# Net content changes:
#   1) The "I want my vote in the poll to be anonymous." user story switches
#      from yellow highlight to bright-green highlight (done -> matches the
#      other already-completed "green" stories around it).
#   2) Word's "last edit position" bookmark (_GoBack) moves from the end of
#      the Google-Maps user story paragraph to the middle of the word
#      "want" inside the final "record of who has seen my announcements"
#      user story paragraph (an artifact of where the author's cursor was
#      when the document was last saved).

$d = $word.ActiveDocument

# --- 1) Highlight: yellow -> bright green -----------------------------
$rng = $d.Content
$target = "As a roommate or household administrator, I want my vote in the poll to be anonymous."
$null = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Font.HighlightColorIndex = 4   # wdBrightGreen -> w:highlight w:val="green"
}

# --- 2) Relocate the _GoBack bookmark ----------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng2 = $d.Content
$null = $rng2.Find.Execute("As an announcement poster, I wan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $collapsed = $d.Range($rng2.End, $rng2.End)
    $d.Bookmarks.Add("_GoBack", $collapsed)
}

Write-Output "edit complete"
